# pred_par.xlsx — add SVR parameter loading columns (svr_kernel_scale,
# svr_epsilon, svr_box_constraint) next to the existing RNN parameter
# table, and tidy up a few cells whose explicit "default" style is no
# longer needed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New SVR parameter columns (K:M), mirroring the existing header/value
#     row-1/row-2 layout used for the RNN parameters in columns A:J. ---
$ws.Range("K1").Value = "svr_kernel_scale"
$ws.Range("L1").Value = "svr_epsilon"
$ws.Range("M1").Value = "svr_box_constraint"

$ws.Range("K2").Value = 100
$ws.Range("L2").Value = 0.1
$ws.Range("M2").Value = 5

# --- Drop the redundant explicit "default" cell style from a handful of
#     cells that only ever carried it as a no-op formatting choice. ---
$ws.Range("B1").Style = "Normal"
$ws.Range("C1").Style = "Normal"
$ws.Range("H1").Style = "Normal"
$ws.Range("H2").Style = "Normal"
$ws.Range("A5").Style = "Normal"
$ws.Range("A6").Style = "Normal"
$ws.Range("A7").Style = "Normal"
$ws.Range("A8").Style = "Normal"
$ws.Range("A11").Style = "Normal"
$ws.Range("D14:G16").Style = "Normal"
$ws.Range("I14:I16").Style = "Normal"

# Row 7 no longer needs its explicit row-level custom formatting either.
$ws.Rows.Item(7).ClearFormats()

# A13 was just an empty, styled placeholder cell — clear it out entirely.
$ws.Range("A13").Style = "Normal"
$ws.Range("A13").ClearContents()

# Selection cursor ends up parked on H10.
$ws.Range("H10").Select()
